$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "Per implementare ... algebra relazionale: " ->
#    split the trailing run into three runs, removing the ": " suffix and
#    adding a separate "." run plus a separate trailing " " run.
# ---------------------------------------------------------------------------
$paraA = $d.Content
$paraA.Find.Execute("Per implementare", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraA.Expand(4)  # wdParagraph

$paraA_xmlFrag = '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>Per implementare le interrogazioni richieste dalla traccia del progetto, sono state utilizzate varie espressioni,</w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> di seguito riportate con l’equivalente espressione scritta in algebra relazionale</w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t>.</w:t></w:r><w:r><w:rPr><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p>'
$paraA_xmlWrap = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraA_xmlFrag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraA.InsertXML($paraA_xmlWrap)

# ---------------------------------------------------------------------------
# 2) Paragraph "Query 1:" (a ListParagraph numbered item) is replaced by:
#      - an empty plain paragraph
#      - a bold "Query 1" paragraph (no colon)
#      - an italic paragraph with the SELECT clause
#      - an italic paragraph with the FROM clause
#      - an italic paragraph with the WHERE clause
#      - an empty italic paragraph
# ---------------------------------------------------------------------------
$paraB = $d.Content
$paraB.Find.Execute("Query 1:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraB.Expand(4)  # wdParagraph

$paraB_xmlFrag = '<w:p><w:pPr><w:rPr><w:lang w:val="it-IT"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Query 1</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.CodLibro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.Titolo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.ISBN</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.Lingua</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.AnnoPubb</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>Libro.CodDip</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve">FROM </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>BibliotecaUNIFE.Libro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>WHERE Titolo LIKE ''%</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>" .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> $</w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t>nomeLibro</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> .</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr><w:t xml:space="preserve"> "%''";</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/><w:iCs/><w:lang w:val="it-IT"/></w:rPr></w:pPr></w:p>'
$paraB_xmlWrap = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $paraB_xmlFrag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$paraB.InsertXML($paraB_xmlWrap)

Write-Host "Done applying edits"
